# "CO cpu fix problem about "1""
# The combinational-logic decode sheet ("组合逻辑译码表") has several
# condition cells that were missing an "&W1" qualifier. This appends
# "&W1" to the existing ST/!ST/!C/!Z conditions that needed it, and
# fills in two previously-blank cells with "ST&W1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("组合逻辑译码表")

$ws.Range("I2").Value  = "!C&W1"
$ws.Range("J2").Value  = "!Z&W1"

$ws.Range("Q2").Value  = "!ST&W1"
$ws.Range("Q4").Value  = "!ST&W1"
$ws.Range("Q5").Value  = "!ST&W1"
$ws.Range("Q11").Value = "!ST&W1"
$ws.Range("Q15").Value = "ST&W1"
$ws.Range("Q18").Value = "!ST&W1"
$ws.Range("Q25").Value = "ST&W1"

$ws.Range("T4").Value  = "!ST&W1"
$ws.Range("T11").Value = "!ST&W1"
$ws.Range("T12").Value = "ST&W1"
$ws.Range("T14").Value = "!ST&W1"
$ws.Range("T16").Value = "ST&W1"

$ws.Range("U4").Value  = "!ST&W1"
$ws.Range("U13").Value = "ST&W1"
$ws.Range("U14").Value = "!ST&W1"
$ws.Range("U16").Value = "ST&W1"
